$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "عنوان السكن" column header to "العنوان" first.
$ws.Range("F1").Value = "العنوان"

# Clear the data rows for the 2nd and 3rd contacts (محمود... and محمد...),
# leaving only the first contact (عمار قصاب) plus blank rows.
$ws.Range("A3:G4").ClearContents()

# Drop the "الرقم الذاتي" column (last table column / column G) entirely,
# shrinking the table from A1:G5 to A1:F5.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item(7).Delete()
$ws.Columns.Item(7).Delete()

# Update the view: no frozen/top-left offset, bigger zoom, new selection.
$ws.Range("E3").Select()
$win = $wb.Windows.Item(1)
$win.Zoom = 210
